# Update KOSS yearly financials sheet with latest figures
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("KOSS")

# Income Statement section
$ws.Range("H20").Value = 0        # Total Other Income/Expenses Net
$ws.Range("H21").Value = -8900    # Earnings Before Interest And Taxes
$ws.Range("H22").Value = "NA"     # Interest Expense
$ws.Range("H32").Value = 0        # Other Items

# Balance Sheet section
$ws.Range("D48").Value = 4200     # Property Plant and Equipment
$ws.Range("D54").Value = 22800    # Total Assets
$ws.Range("D59").Value = 1700     # Other Current Liabilities
$ws.Range("D60").Value = 3200     # Total Current Liabilities
$ws.Range("D62").Value = 5600     # Other Liabilities
$ws.Range("D66").Value = 8600     # Total Liabilities
$ws.Range("D72").Value = 8400     # Retained Earnings
$ws.Range("D76").Value = 14200    # Total Stockholder Equity

# Cash Flow Statement section
$ws.Range("J91").Value = -400     # Capital Expenditures
